$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '60.261.95'
Set-TextCell $ws 'E2' '  -2.95%  '
Set-TextCell $ws 'D3' '3.297.51'
Set-TextCell $ws 'E3' '  -3.62%  '
Set-TextCell $ws 'E4' '  +0.03%  '
Set-TextCell $ws 'D5' '557.17'
Set-TextCell $ws 'E5' '  -3.87%  '
Set-TextCell $ws 'D6' '140.84'
Set-TextCell $ws 'E6' '  -8.39%  '
Set-TextCell $ws 'D8' '3.298.29'
Set-TextCell $ws 'E8' '  -3.56%  '
Set-TextCell $ws 'E9' '  -3.64%  '
Set-TextCell $ws 'D10' '7.88'
Set-TextCell $ws 'E10' '  -2.65%  '
Set-TextCell $ws 'D11' '0.119'
Set-TextCell $ws 'E11' '  -5.18%  '
Set-TextCell $ws 'D12' '0.407'
Set-TextCell $ws 'E12' '  -2.66%  '
Set-TextCell $ws 'D13' '3.860.97'
Set-TextCell $ws 'E14' '  -0.50%  '
Set-TextCell $ws 'D15' '26.50'
Set-TextCell $ws 'E15' '  -7.46%  '
Set-TextCell $ws 'D16' '3.301.98'
Set-TextCell $ws 'E16' '  -3.90%  '
Set-TextCell $ws 'E17' '  -4.92%  '
Set-TextCell $ws 'D18' '60.238.33'
Set-TextCell $ws 'E18' '  -3.02%  '
Set-TextCell $ws 'D19' '6.05'
Set-TextCell $ws 'E19' '  -7.15%  '
Set-TextCell $ws 'D20' '13.63'
Set-TextCell $ws 'E20' '  -5.44%  '
Set-TextCell $ws 'D21' '8.52'
Set-TextCell $ws 'E21' '  -5.24%  '
Set-TextCell $ws 'D22' '373.40'
Set-TextCell $ws 'E22' '  -2.36%  '
Set-TextCell $ws 'E23' '  -0.04%  '
Set-TextCell $ws 'E24' '  -5.18%  '
Set-TextCell $ws 'E25' '  -6.96%  '
Set-TextCell $ws 'D26' '3.434.33'
Set-TextCell $ws 'E26' '  -3.63%  '
Set-TextCell $ws 'E27' '  -9.34%  '
Set-TextCell $ws 'E28' '  -2.26%  '
Set-TextCell $ws 'D29' '1.00'
Set-TextCell $ws 'E29' '  +0.20%  '
Set-TextCell $ws 'D30' '7.05'
Set-TextCell $ws 'E30' '  -7.84%  '
Set-TextCell $ws 'E31' '  +0.00%  '
Set-TextCell $ws 'E32' '  -4.96%  '
Set-TextCell $ws 'E33' '  -6.92%  '
Set-TextCell $ws 'D34' '22.57'
Set-TextCell $ws 'E34' '  -3.05%  '
Set-TextCell $ws 'E35' '  -5.15%  '
Set-TextCell $ws 'B36' 'NEARProtocol'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D36' '5.03'
Set-TextCell $ws 'E36' '  -8.70%  '
Set-TextCell $ws 'B37' 'Monero'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D37' '165.57'
Set-TextCell $ws 'E37' '  -1.86%  '
Set-TextCell $ws 'E38' '  -4.56%  '
Set-TextCell $ws 'D39' '6.63'
Set-TextCell $ws 'E39' '  -4.72%  '
Set-TextCell $ws 'D40' '3.329.62'
Set-TextCell $ws 'E40' '  -3.70%  '
Set-TextCell $ws 'D41' '0.0722'
Set-TextCell $ws 'E41' '  -7.60%  '
Set-TextCell $ws 'D42' '25.49'
Set-TextCell $ws 'E42' '  -17.66%  '
Set-TextCell $ws 'D43' '41.75'
Set-TextCell $ws 'E43' '  -2.42%  '
Set-TextCell $ws 'E44' '  -4.51%  '
Set-TextCell $ws 'E45' '  -4.19%  '
Set-TextCell $ws 'D46' '4.10'
Set-TextCell $ws 'E46' '  -7.14%  '
Set-TextCell $ws 'D47' '1.57'
Set-TextCell $ws 'E47' '  -6.59%  '
Set-TextCell $ws 'E48' '  -0.02%  '
Set-TextCell $ws 'D49' '2.319.64'
Set-TextCell $ws 'E49' '  -9.39%  '
Set-TextCell $ws 'E50' '  -6.26%  '
Set-TextCell $ws 'D51' '21.47'
Set-TextCell $ws 'E51' '  -7.62%  '
